# Adds a new "2022-Q1" fund-holdings detail sheet (taking over the slot that
# used to be the old "总计" sheet) and recreates the "总计" (grand-total)
# sheet at the end of the workbook with a new leading "2022-Q1" row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Step 1: duplicate the existing "总计" sheet so we keep its ready-made
# layout/formatting, then turn the duplicate into the refreshed totals
# table (old rows shifted down one, new 2022-Q1 row inserted at the top).
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Copy([System.Reflection.Missing]::Value, $totalSheet)
$newTotal = $wb.Worksheets.Item("总计 (2)")
$newTotal.Name = "总计TEMP"

for ($r = 6; $r -ge 2; $r--) {
  $dest = $r + 1
  $newTotal.Range("A$dest").Value = $newTotal.Range("A$r").Value()
  $newTotal.Range("B$dest").Value = $newTotal.Range("B$r").Value()
  $newTotal.Range("C$dest").Value = $newTotal.Range("C$r").Value()
  $newTotal.Range("D$dest").Value = $newTotal.Range("D$r").Value()
}

$newTotal.Range("A2").Value = 0
$newTotal.Range("B2").Value = "2022-Q1"
$newTotal.Range("C2").Value = 8
$newTotal.Range("D2").Value = 6.39

$newTotal.Range("A2").Copy()
$newTotal.Range("A3:A7").PasteSpecial(-4122)  # xlPasteFormats
for ($r = 2; $r -le 7; $r++) {
  $newTotal.Range("A$r").Value = $r - 2
}

# ---------------------------------------------------------------------------
# Step 2: repurpose the original "总计" sheet as the new "2022-Q1" detail
# sheet (same layout used by the other quarterly fund-holdings sheets).
# ---------------------------------------------------------------------------
$totalSheet.Name = "2022-Q1"
$totalSheet.Cells.Clear()

$src = $wb.Worksheets.Item("2021-Q4")
$src.Range("B1:H1").Copy()
$totalSheet.Range("B1").PasteSpecial(-4122)
$src.Range("A2:A9").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)

$totalSheet.Range("B1").Value = "基金代码"
$totalSheet.Range("C1").Value = "基金名称"
$totalSheet.Range("D1").Value = "基金规模"
$totalSheet.Range("E1").Value = "股票总仓位"
$totalSheet.Range("F1").Value = "仓位占比"
$totalSheet.Range("G1").Value = "持有市值(亿元)"
$totalSheet.Range("H1").Value = "仓位排名"

$data = @(
  @("001645", "国泰大健康股票A", "34.81", "90.83", "6.93", "2.4123", 5),
  @("020001", "国泰金鹰增长灵活配置混合", "17.38", "91.77", "8.17", "1.4199", 5),
  @("009805", "国泰医药健康股票A", "12.14", "92.49", "8.56", "1.0392", 5),
  @("160215", "国泰价值经典灵活配置混合（LOF）", "6.36", "92.86", "7.29", "0.4636", 7),
  @("008370", "国泰研究精选两年持有期混合", "4.15", "92.87", "8.88", "0.3685", 6),
  @("009804", "国泰研究优势混合", "4.14", "90.91", "8.52", "0.3527", 5),
  @("011321", "国泰大健康股票C", "3.47", "90.83", "6.93", "0.2405", 5),
  @("011326", "国泰医药健康股票C", "1.09", "92.49", "8.56", "0.0933", 5)
)

$row = 2
foreach ($item in $data) {
  $totalSheet.Cells.Item($row, 1).Value = $row - 2
  $totalSheet.Cells.Item($row, 2).Value = "'" + $item[0]
  $totalSheet.Cells.Item($row, 3).Value = $item[1]
  $totalSheet.Cells.Item($row, 4).Value = "'" + $item[2]
  $totalSheet.Cells.Item($row, 5).Value = "'" + $item[3]
  $totalSheet.Cells.Item($row, 6).Value = "'" + $item[4]
  $totalSheet.Cells.Item($row, 7).Value = "'" + $item[5]
  $totalSheet.Cells.Item($row, 8).Value = $item[6]
  $row++
}

# ---------------------------------------------------------------------------
# Step 3: the finished totals duplicate becomes the "总计" sheet again, and
# we restore the original active sheet/tab selection.
# ---------------------------------------------------------------------------
$newTotal.Name = "总计"

$wb.Worksheets.Item(1).Activate()
